$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 115
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H69").Value = 4200
$ws.Range("I69").Value = 2800
$ws.Range("K69").Value = 8400
$ws.Range("M69").Value = -7526

$ws.Range("H70").Value = 12220.444
$ws.Range("I70").Value = 1358
$ws.Range("K70").Value = 4074
$ws.Range("M70").Value = -3804

$ws.Range("H72").Value = 4200
$ws.Range("I72").Value = 2800
$ws.Range("K72").Value = 25200
$ws.Range("M72").Value = -20832

$ws.Range("H73").Value = 12220.444
$ws.Range("I73").Value = 1358
$ws.Range("K73").Value = 4074
$ws.Range("M73").Value = -3138

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H94").Value = 3100.1428
$ws.Range("I94").Value = 2783.5
$ws.Range("K94").Value = 2783.5
$ws.Range("M94").Value = -2332.5

$ws.Range("H129").Value = 852.9818
$ws.Range("J129").Value = 871.9388
$ws.Range("L129").Value = 2615.8164
$ws.Range("N129").Value = -12615.8164

$ws.Range("H137").Value = 43200
$ws.Range("I137").Value = 1450.4445
$ws.Range("J137").Value = 168448.67
$ws.Range("K137").Value = 4351.333500000001
$ws.Range("L137").Value = 505346.01
$ws.Range("M137").Value = -1801.333500000001
$ws.Range("N137").Value = -510446.01

$ws.Range("H138").Value = 3264.8
$ws.Range("I138").Value = 2936.9048
$ws.Range("J138").Value = 4029.889
$ws.Range("K138").Value = 8810.714399999999
$ws.Range("L138").Value = 12089.667
$ws.Range("M138").Value = -3670.714399999999
$ws.Range("N138").Value = -22369.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1938943.2
$ws.Range("I2").Value = 2907840
$ws.Range("K2").Value = 2907840
$ws.Range("M2").Value = -2907727

$ws.Range("H32").Value = 3172.9275
$ws.Range("I32").Value = 2429
$ws.Range("K32").Value = 2429
$ws.Range("M32").Value = -2142

$ws.Range("H116").Value = 1938943.2
$ws.Range("I116").Value = 2907840
$ws.Range("K116").Value = 2907840
$ws.Range("M116").Value = -2905546

$ws.Range("H132").Value = 1635.1724
$ws.Range("I132").Value = 1082.4762
$ws.Range("K132").Value = 3247.4286
$ws.Range("M132").Value = -717.4286000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1938943.2
$ws.Range("I3").Value = 2907840
$ws.Range("K3").Value = 2907840
$ws.Range("M3").Value = -2907726

$ws.Range("H86").Value = 170033.67
$ws.Range("I86").Value = 2533.3333
$ws.Range("J86").Value = 225867.11
$ws.Range("K86").Value = 2533.3333
$ws.Range("L86").Value = 225867.11
$ws.Range("M86").Value = -1410.3333
$ws.Range("N86").Value = -228113.11

$ws.Range("H89").Value = 170033.67
$ws.Range("I89").Value = 2533.3333
$ws.Range("J89").Value = 225867.11
$ws.Range("K89").Value = 12666.6665
$ws.Range("L89").Value = 1129335.55
$ws.Range("M89").Value = -7050.666499999999
$ws.Range("N89").Value = -1140567.55

$ws.Range("H107").Value = 2108.8125
$ws.Range("J107").Value = 2200
$ws.Range("L107").Value = 2200
$ws.Range("N107").Value = -6040

$ws.Range("H134").Value = 5415.9395
$ws.Range("I134").Value = 5890.552
$ws.Range("J134").Value = 1975
$ws.Range("K134").Value = 17671.656
$ws.Range("L134").Value = 5925
$ws.Range("M134").Value = -15136.656
$ws.Range("N134").Value = -10995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2499.375
$ws.Range("I31").Value = 1900
$ws.Range("J31").Value = 2699.1667
$ws.Range("K31").Value = 1900
$ws.Range("L31").Value = 2699.1667
$ws.Range("M31").Value = -1605
$ws.Range("N31").Value = -3289.1667

$ws.Range("H34").Value = 2499.375
$ws.Range("I34").Value = 1900
$ws.Range("J34").Value = 2699.1667
$ws.Range("K34").Value = 1900
$ws.Range("L34").Value = 2699.1667
$ws.Range("M34").Value = -1698
$ws.Range("N34").Value = -3103.1667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 618.625
$ws.Range("I26").Value = 723.75
$ws.Range("K26").Value = 2171.25
$ws.Range("M26").Value = -1883.25

$ws.Range("H104").Value = 2996.652
$ws.Range("J104").Value = 3228.9
$ws.Range("L104").Value = 9686.7
$ws.Range("N104").Value = -14928.7

$ws.Range("H107").Value = 700.86957
$ws.Range("J107").Value = 849.05884
$ws.Range("L107").Value = 2547.17652
$ws.Range("N107").Value = -6387.17652

$ws.Range("H131").Value = 818.6531
$ws.Range("J131").Value = 818.6531
$ws.Range("L131").Value = 2455.9593
$ws.Range("N131").Value = -12535.9593

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2352.182
$ws.Range("J46").Value = 3195.6667
$ws.Range("L46").Value = 3195.6667
$ws.Range("N46").Value = -3571.6667

$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H82").Value = 2117.8667
$ws.Range("I82").Value = 1556.2858
$ws.Range("J82").Value = 2609.25
$ws.Range("K82").Value = 1556.2858
$ws.Range("L82").Value = 2609.25
$ws.Range("M82").Value = -1195.2858
$ws.Range("N82").Value = -3331.25

$ws.Range("H85").Value = 2117.8667
$ws.Range("I85").Value = 1556.2858
$ws.Range("J85").Value = 2609.25
$ws.Range("K85").Value = 1556.2858
$ws.Range("L85").Value = 2609.25
$ws.Range("M85").Value = -308.2858000000001
$ws.Range("N85").Value = -5105.25

$ws.Range("H136").Value = 2524.48
$ws.Range("I136").Value = 1534
$ws.Range("K136").Value = 4602
$ws.Range("M136").Value = -2052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2680
$ws.Range("J96").Value = 2600
$ws.Range("L96").Value = 2600
$ws.Range("N96").Value = -5346

$ws.Range("H107").Value = 827.375
$ws.Range("I107").Value = 710.38464
$ws.Range("K107").Value = 2131.15392
$ws.Range("M107").Value = -211.1539199999997
